$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.323.97'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.834.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.09%  '
$ws.Range('E4').Value = '  +1.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4744'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.80%  '
$ws.Range('E8').Value = '  +0.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07451'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8855'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.92%  '
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.875.92'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07334'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.444'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.65'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.578'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.18%  '
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008797'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.464.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.291'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('E23').Value = '  +0.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.085.62'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.42%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.02'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('E27').Value = '  +1.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.147'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.235'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08995'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.177'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7499'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.546'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.945'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.27%  '
$ws.Range('E36').Value = '  +0.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.102'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05347'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.35%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01955'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.964'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.261'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.380'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5313'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.69%  '
$ws.Range('E44').Value = '  +0.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.486'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4921'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.57'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '105.08'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.86%  '
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.672'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06293'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.10%  '
